# Update countries & provincias Spain
# Applies the 18-May-2020 14:35 data refresh to the "Pais" sheet:
#  - Updated case counts for several countries (rows 8, 11, 24, 25, 48,
#    54-57, 60, 82, 144)
#  - Catar overtakes Ecuador in the ranking (rows 24/25 swap country names
#    while each keeps its own updated figures)
#  - Barein overtakes Afganistan/Australia/Argelia, pushing them down one
#    row each (rows 54-57 cascade)
#  - Footer timestamp bumped from 14:05 to 14:35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 14:35"

# Row 8 (Brasil)
$ws.Range("B8").Value = 242313
$ws.Range("C8").Value = 1233
$ws.Range("E8").Value = 131999
$ws.Range("G8").Value = 74
$ws.Range("H8").Value = 16192

# Row 11 (Alemania)
$ws.Range("B11").Value = 176788
$ws.Range("C11").Value = 137
$ws.Range("E11").Value = 14136
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 8052

# Rows 24/25: Catar overtakes Ecuador
$ws.Range("A24").Value = "Catar"
$ws.Range("B24").Value = 33969
$ws.Range("C24").Value = 1365
$ws.Range("D24").Value = 4899
$ws.Range("E24").Value = 29055
$ws.Range("H24").Value = 15

$ws.Range("A25").Value = "Ecuador"
$ws.Range("B25").Value = 33182
$ws.Range("D25").Value = 3433
$ws.Range("E25").Value = 27013
$ws.Range("H25").Value = 2736

# Row 48 (Dinamarca)
$ws.Range("D48").Value = 9301
$ws.Range("E48").Value = 1119
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 548

# Rows 54-57: Barein overtakes Afganistan / Australia / Argelia
$ws.Range("A54").Value = "Barein"
$ws.Range("B54").Value = 7156
$ws.Range("C54").Value = 200
$ws.Range("D54").Value = 2929
$ws.Range("E54").Value = 4215
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 12

$ws.Range("A55").Value = "Afganistan"
$ws.Range("B55").Value = 7072
$ws.Range("C55").Value = 408
$ws.Range("D55").Value = 801
$ws.Range("E55").Value = 6098
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 173

$ws.Range("A56").Value = "Australia"
$ws.Range("B56").Value = 7060
$ws.Range("C56").Value = 15
$ws.Range("D56").Value = 6392
$ws.Range("E56").Value = 569
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 99

$ws.Range("A57").Value = "Argelia"
$ws.Range("B57").Value = 7019
$ws.Range("D57").Value = 3507
$ws.Range("E57").Value = 2964
$ws.Range("H57").Value = 548

# Row 60 (Kazajistan)
$ws.Range("D60").Value = 3373
$ws.Range("E60").Value = 3033

# Row 82 (Croacia)
$ws.Range("B82").Value = 2228
$ws.Range("C82").Value = 2
$ws.Range("D82").Value = 1946
$ws.Range("E82").Value = 187

# Row 144 (Madagascar)
$ws.Range("B144").Value = 322
$ws.Range("C144").Value = 18
$ws.Range("D144").Value = 119
$ws.Range("E144").Value = 202
